# Updated symbol list on Tue Jan 24 21:08:15 UTC 2023 with GitHub Actions
# Applies the refreshed coinranking.com price/volume/rank snapshot to Sheet1.
#
# Values are written as plain text (not numbers) to match the workbook's
# existing inline-string cell format: setting NumberFormat to "@" (Text)
# before assigning .Value keeps Excel from auto-coercing numeric-looking
# strings (prices, "NN" rank values, "N.NN%" percentages) into real
# numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "310.55" },
    @{ Cell = "E2"; Value = "1.65%" },
    @{ Cell = "G2"; Value = "21" },
    @{ Cell = "E3"; Value = "-2.10%" },
    @{ Cell = "G3"; Value = "21" },
    @{ Cell = "D4"; Value = "5.121" },
    @{ Cell = "E4"; Value = "1.49%" },
    @{ Cell = "G4"; Value = "21" },
    @{ Cell = "D5"; Value = "0.08205" },
    @{ Cell = "E5"; Value = "3.50%" },
    @{ Cell = "G5"; Value = "21" },
    @{ Cell = "D6"; Value = "2.028" },
    @{ Cell = "E6"; Value = "-9.55%" },
    @{ Cell = "G6"; Value = "21" },
    @{ Cell = "D7"; Value = "7.977" },
    @{ Cell = "E7"; Value = "-0.16%" },
    @{ Cell = "G7"; Value = "21" },
    @{ Cell = "B8"; Value = "BTSEToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D8"; Value = "2.949" },
    @{ Cell = "E8"; Value = "12.00%" },
    @{ Cell = "G8"; Value = "21" },
    @{ Cell = "B9"; Value = "MXToken" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D9"; Value = "0.9283" },
    @{ Cell = "E9"; Value = "0.06%" },
    @{ Cell = "G9"; Value = "21" },
    @{ Cell = "B10"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D10"; Value = "0.1081" },
    @{ Cell = "E10"; Value = "9.91%" },
    @{ Cell = "G10"; Value = "21" },
    @{ Cell = "B11"; Value = "WazirX" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D11"; Value = "0.1925" },
    @{ Cell = "E11"; Value = "2.71%" },
    @{ Cell = "G11"; Value = "21" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.09418" },
    @{ Cell = "E12"; Value = "3.96%" },
    @{ Cell = "G12"; Value = "21" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D13"; Value = "0.03588" },
    @{ Cell = "E13"; Value = "-4.30%" },
    @{ Cell = "G13"; Value = "21" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D14"; Value = "0.09909" },
    @{ Cell = "E14"; Value = "-0.32%" },
    @{ Cell = "G14"; Value = "21" },
    @{ Cell = "B15"; Value = "BitForexToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D15"; Value = "0.001431" },
    @{ Cell = "E15"; Value = "-1.12%" },
    @{ Cell = "G15"; Value = "21" },
    @{ Cell = "B16"; Value = "TigerCash" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D16"; Value = "0.005760" },
    @{ Cell = "E16"; Value = "1.94%" },
    @{ Cell = "G16"; Value = "21" },
    @{ Cell = "B17"; Value = "LEO" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D17"; Value = "3.474" },
    @{ Cell = "E17"; Value = "0.25%" },
    @{ Cell = "G17"; Value = "21" },
    @{ Cell = "B18"; Value = "GateToken" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D18"; Value = "4.128" },
    @{ Cell = "E18"; Value = "-0.42%" },
    @{ Cell = "G18"; Value = "21" },
    @{ Cell = "D19"; Value = "0.3418" },
    @{ Cell = "E19"; Value = "1.43%" },
    @{ Cell = "G19"; Value = "21" },
    @{ Cell = "D20"; Value = "0.1320" },
    @{ Cell = "E20"; Value = "0.27%" },
    @{ Cell = "G20"; Value = "21" },
    @{ Cell = "D21"; Value = "5.091" },
    @{ Cell = "E21"; Value = "0.34%" },
    @{ Cell = "G21"; Value = "21" },
    @{ Cell = "D22"; Value = "0.2191" },
    @{ Cell = "E22"; Value = "-2.59%" },
    @{ Cell = "G22"; Value = "21" },
    @{ Cell = "D23"; Value = "0.04548" },
    @{ Cell = "E23"; Value = "-0.58%" },
    @{ Cell = "G23"; Value = "21" },
    @{ Cell = "E24"; Value = "-0.59%" },
    @{ Cell = "G24"; Value = "21" },
    @{ Cell = "D25"; Value = "0.004788" },
    @{ Cell = "E25"; Value = "0.19%" },
    @{ Cell = "G25"; Value = "21" },
    @{ Cell = "D26"; Value = "0.0001249" },
    @{ Cell = "E26"; Value = "-3.96%" },
    @{ Cell = "G26"; Value = "21" },
    @{ Cell = "D27"; Value = "0.0004449" },
    @{ Cell = "E27"; Value = "-6.11%" },
    @{ Cell = "G27"; Value = "21" },
    @{ Cell = "G28"; Value = "21" },
    @{ Cell = "G29"; Value = "21" },
    @{ Cell = "G30"; Value = "21" },
    @{ Cell = "G31"; Value = "21" },
    @{ Cell = "G32"; Value = "21" },
    @{ Cell = "G33"; Value = "21" },
    @{ Cell = "G34"; Value = "21" },
    @{ Cell = "G35"; Value = "21" },
    @{ Cell = "G36"; Value = "21" },
    @{ Cell = "G37"; Value = "21" },
    @{ Cell = "G38"; Value = "21" },
    @{ Cell = "E39"; Value = "3.05%" },
    @{ Cell = "G39"; Value = "21" },
    @{ Cell = "D40"; Value = "0.04908" },
    @{ Cell = "E40"; Value = "0.07%" },
    @{ Cell = "G40"; Value = "21" },
    @{ Cell = "D41"; Value = "0.007828" },
    @{ Cell = "E41"; Value = "0.41%" },
    @{ Cell = "G41"; Value = "21" },
    @{ Cell = "D42"; Value = "0.009867" },
    @{ Cell = "E42"; Value = "26.56%" },
    @{ Cell = "G42"; Value = "21" },
    @{ Cell = "D43"; Value = "0.1383" },
    @{ Cell = "E43"; Value = "-0.82%" },
    @{ Cell = "G43"; Value = "21" },
    @{ Cell = "D44"; Value = "0.002114" },
    @{ Cell = "E44"; Value = "-0.94%" },
    @{ Cell = "G44"; Value = "21" },
    @{ Cell = "D45"; Value = "0.01156" },
    @{ Cell = "E45"; Value = "1.02%" },
    @{ Cell = "G45"; Value = "21" },
    @{ Cell = "D46"; Value = "0.00006480" },
    @{ Cell = "E46"; Value = "5.08%" },
    @{ Cell = "G46"; Value = "21" },
    @{ Cell = "D47"; Value = "0.00000000749" },
    @{ Cell = "E47"; Value = "-0.10%" },
    @{ Cell = "G47"; Value = "21" },
    @{ Cell = "D48"; Value = "64.72" },
    @{ Cell = "E48"; Value = "25.03%" },
    @{ Cell = "G48"; Value = "21" },
    @{ Cell = "D49"; Value = "0.001300" },
    @{ Cell = "E49"; Value = "-27.79%" },
    @{ Cell = "G49"; Value = "21" },
    @{ Cell = "D50"; Value = "0.00002099" },
    @{ Cell = "E50"; Value = "-0.10%" },
    @{ Cell = "G50"; Value = "21" },
    @{ Cell = "D51"; Value = "0.0001999" },
    @{ Cell = "E51"; Value = "-0.10%" },
    @{ Cell = "G51"; Value = "21" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.NumberFormat = "@"
    $range.Value = $u.Value
}

Write-Host ("Updated {0} cells" -f $updates.Count)
